# Update odds/score values in Sheet1 to reflect the latest FlashScore scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    "G3"  = 2.8
    "I3"  = 2.3
    "L3"  = 3
    "AJ3" = 6.5
    "AN3" = 13
    "AP3" = 23

    "G4"  = 6.5
    "H4"  = 4.75
    "I4"  = 1.45
    "J4"  = 6
    "K4"  = 2.5
    "L4"  = 1.95
    "N4"  = 17
    "O4"  = 1.17
    "P4"  = 5
    "S4"  = 1.53
    "T4"  = 2.4
    "U4"  = 1.88
    "V4"  = 1.98
    "W4"  = 2.25
    "X4"  = 1.57
    "Y4"  = 1.29
    "Z4"  = 3.5
    "AA4" = 1.7
    "AB4" = 2.05
    "AC4" = 21
    "AD4" = 34
    "AE4" = 19
    "AF4" = 67
    "AG4" = 41
    "AH4" = 41
    "AI4" = 17
    "AJ4" = 9
    "AK4" = 15
    "AN4" = 8
    "AP4" = 11
    "AQ4" = 11

    "G5"  = 2.52
    "H5"  = 3.35
    "I5"  = 2.52
    "J5"  = 3.05
    "K5"  = 2.15
    "L5"  = 3.1
    "O5"  = 1.23
    "P5"  = 3.4
    "S5"  = 1.7
    "T5"  = 1.93
    "W5"  = 2.57
    "X5"  = 1.38
    "AA5" = 1.55
    "AB5" = 2.15
    "AC5" = 10.25
    "AF5" = 29
    "AG5" = 19
    "AH5" = 24
    "AI5" = 11.75
    "AJ5" = 6.7
    "AK5" = 12
    "AL5" = 45
    "AM5" = 9.75
    "AN5" = 13.5
    "AQ5" = 20
    "AR5" = 26
    "AS5" = 300

    "G7"  = 2.3
    "H7"  = 3.2
    "I7"  = 2.88
    "J7"  = 3
    "L7"  = 3.5
    "M7"  = 1.04
    "N7"  = 9
    "O7"  = 1.3
    "P7"  = 3.4
    "S7"  = 2.03
    "T7"  = 1.78
    "AA7" = 1.8
    "AB7" = 1.91
    "AC7" = 8
    "AD7" = 11
    "AE7" = 10
    "AF7" = 21
    "AN7" = 15
    "AP7" = 29
    "AQ7" = 23
    "AR7" = 34
}

foreach ($ref in $updates.Keys) {
    $ws.Range($ref).Value = $updates[$ref]
}
